# Auto-generated edit script: apply 2022-08-07 data update across workbook
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("Citywide Totals")
$ws.Range("I2").Value = 4232
$ws.Range("I3").Value = 4444
$ws.Range("B4").Value = 1656
$ws.Range("H4").Value = 1669
$ws.Range("I4").Value = 1024
$ws.Range("I5").Value = 406
$ws.Range("I6").Value = 4846
$ws.Range("B7").Value = 23288
$ws.Range("H7").Value = 25980
$ws.Range("I7").Value = 14952

$ws = $wb.Worksheets.Item("Uptown")
$ws.Range("I6").Value = 60
$ws.Range("I7").Value = 169

$ws = $wb.Worksheets.Item("West Ridge")
$ws.Range("I6").Value = 58
$ws.Range("I7").Value = 156

$ws = $wb.Worksheets.Item("Bridgeport")
$ws.Range("I3").Value = 23
$ws.Range("I7").Value = 82

$ws = $wb.Worksheets.Item("Grand Crossing")
$ws.Range("I2").Value = 154
$ws.Range("I3").Value = 158
$ws.Range("I7").Value = 489

$ws = $wb.Worksheets.Item("Woodlawn")
$ws.Range("I4").Value = 21
$ws.Range("I7").Value = 286

$ws = $wb.Worksheets.Item("North Lawndale")
$ws.Range("I2").Value = 143
$ws.Range("I6").Value = 188
$ws.Range("I7").Value = 589

$ws = $wb.Worksheets.Item("Gage Park")
$ws.Range("I6").Value = 51
$ws.Range("I7").Value = 146

$ws = $wb.Worksheets.Item("By Neighborhood")
$ws.Range("I6").Value = 95
$ws.Range("I7").Value = 475
$ws.Range("I8").Value = 904
$ws.Range("I9").Value = 68
$ws.Range("I10").Value = 102
$ws.Range("I11").Value = 226
$ws.Range("I14").Value = 82
$ws.Range("I19").Value = 413
$ws.Range("I24").Value = 39
$ws.Range("I25").Value = 74
$ws.Range("B27").Value = 110
$ws.Range("I27").Value = 135
$ws.Range("I29").Value = 957
$ws.Range("I31").Value = 146
$ws.Range("I33").Value = 691
$ws.Range("I34").Value = 69
$ws.Range("I37").Value = 489
$ws.Range("I39").Value = 12
$ws.Range("I41").Value = 66
$ws.Range("I42").Value = 507
$ws.Range("I46").Value = 30
$ws.Range("I48").Value = 206
$ws.Range("I51").Value = 159
$ws.Range("I52").Value = 324
$ws.Range("I54").Value = 337
$ws.Range("H55").Value = 309
$ws.Range("I60").Value = 70
$ws.Range("I63").Value = 54
$ws.Range("I64").Value = 131
$ws.Range("I67").Value = 589
$ws.Range("I68").Value = 48
$ws.Range("I71").Value = 45
$ws.Range("I73").Value = 123
$ws.Range("I76").Value = 222
$ws.Range("I77").Value = 83
$ws.Range("I79").Value = 408
$ws.Range("I83").Value = 305
$ws.Range("I85").Value = 674
$ws.Range("I86").Value = 88
$ws.Range("I88").Value = 135
$ws.Range("I89").Value = 169
$ws.Range("I90").Value = 182
$ws.Range("I91").Value = 179
$ws.Range("I92").Value = 44
$ws.Range("I93").Value = 89
$ws.Range("I96").Value = 156
$ws.Range("I97").Value = 112
$ws.Range("I99").Value = 286
$ws.Range("I100").Value = 27
$ws.Range("B101").Value = 23288
$ws.Range("H101").Value = 25980
$ws.Range("I101").Value = 14952

$ws = $wb.Worksheets.Item("South Chicago")
$ws.Range("I3").Value = 121
$ws.Range("I6").Value = 56
$ws.Range("I7").Value = 305

$ws = $wb.Worksheets.Item("Garfield Park")
$ws.Range("I2").Value = 159
$ws.Range("I3").Value = 258
$ws.Range("I5").Value = 31
$ws.Range("I6").Value = 211
$ws.Range("I7").Value = 691

$ws = $wb.Worksheets.Item("Loop")
$ws.Range("I2").Value = 76
$ws.Range("I6").Value = 166
$ws.Range("I7").Value = 337

$ws = $wb.Worksheets.Item("Englewood")
$ws.Range("I2").Value = 275
$ws.Range("I3").Value = 334
$ws.Range("I4").Value = 49
$ws.Range("I7").Value = 957

$ws = $wb.Worksheets.Item("Chatham")
$ws.Range("I2").Value = 150
$ws.Range("I6").Value = 116
$ws.Range("I7").Value = 413

$ws = $wb.Worksheets.Item("Lake View")
$ws.Range("I4").Value = 22
$ws.Range("I7").Value = 206

$ws = $wb.Worksheets.Item("River North")
$ws.Range("I6").Value = 94
$ws.Range("I7").Value = 222

$ws = $wb.Worksheets.Item("South Shore")
$ws.Range("I3").Value = 271
$ws.Range("I4").Value = 39
$ws.Range("I6").Value = 168
$ws.Range("I7").Value = 674

$ws = $wb.Worksheets.Item("Ashburn")
$ws.Range("I3").Value = 28
$ws.Range("I6").Value = 20
$ws.Range("I7").Value = 95

$ws = $wb.Worksheets.Item("Hermosa")
$ws.Range("I2").Value = 22
$ws.Range("I7").Value = 66

$ws = $wb.Worksheets.Item("Humboldt Park")
$ws.Range("I2").Value = 136
$ws.Range("I3").Value = 175
$ws.Range("I7").Value = 507

$ws = $wb.Worksheets.Item("Avondale")
$ws.Range("I2").Value = 36
$ws.Range("I7").Value = 102

$ws = $wb.Worksheets.Item("Lower West Side")
$ws.Range("H4").Value = 21
$ws.Range("H7").Value = 309

$ws = $wb.Worksheets.Item("Dunning")
$ws.Range("I6").Value = 6
$ws.Range("I7").Value = 39

$ws = $wb.Worksheets.Item("Jefferson Park")
$ws.Range("I4").Value = 1
$ws.Range("I7").Value = 30

$ws = $wb.Worksheets.Item("Washington Park")
$ws.Range("I2").Value = 56
$ws.Range("I7").Value = 179

$ws = $wb.Worksheets.Item("Roseland")
$ws.Range("I2").Value = 121
$ws.Range("I7").Value = 408

$ws = $wb.Worksheets.Item("Near South Side")
$ws.Range("I3").Value = 39
$ws.Range("I7").Value = 131

$ws = $wb.Worksheets.Item("West Lawn")
$ws.Range("I3").Value = 22
$ws.Range("I7").Value = 89

$ws = $wb.Worksheets.Item("Wrigleyville")
$ws.Range("I2").Value = 6
$ws.Range("I6").Value = 27

$ws = $wb.Worksheets.Item("Little Village")
$ws.Range("I3").Value = 113
$ws.Range("I7").Value = 324

$ws = $wb.Worksheets.Item("Garfield Ridge")
$ws.Range("I3").Value = 24
$ws.Range("I7").Value = 69

$ws = $wb.Worksheets.Item("East Side")
$ws.Range("I6").Value = 23
$ws.Range("I7").Value = 74

$ws = $wb.Worksheets.Item("Greektown")
$ws.Range("I5").Value = 8
$ws.Range("I6").Value = 12

$ws = $wb.Worksheets.Item("Belmont Cragin")
$ws.Range("I2").Value = 100
$ws.Range("I7").Value = 226

$ws = $wb.Worksheets.Item("Avalon Park")
$ws.Range("I2").Value = 27
$ws.Range("I6").Value = 15
$ws.Range("I7").Value = 68

$ws = $wb.Worksheets.Item("Portage Park")
$ws.Range("I3").Value = 37
$ws.Range("I6").Value = 30
$ws.Range("I7").Value = 123

$ws = $wb.Worksheets.Item("West Town")
$ws.Range("I3").Value = 21
$ws.Range("I7").Value = 112

$ws = $wb.Worksheets.Item("West Elsdon")
$ws.Range("I6").Value = 19
$ws.Range("I7").Value = 44

$ws = $wb.Worksheets.Item("United Center")
$ws.Range("I2").Value = 37
$ws.Range("I7").Value = 135

$ws = $wb.Worksheets.Item("Austin")
$ws.Range("I2").Value = 284
$ws.Range("I3").Value = 258
$ws.Range("I6").Value = 285
$ws.Range("I7").Value = 904

$ws = $wb.Worksheets.Item("Edgewater")
$ws.Range("I3").Value = 26
$ws.Range("B4").Value = 14
$ws.Range("I4").Value = 17
$ws.Range("B7").Value = 110
$ws.Range("I7").Value = 135

$ws = $wb.Worksheets.Item("Streeterville")
$ws.Range("I4").Value = 43
$ws.Range("I7").Value = 88

$ws = $wb.Worksheets.Item("Washington Heights")
$ws.Range("I3").Value = 39
$ws.Range("I7").Value = 182

$ws = $wb.Worksheets.Item("Little Italy, UIC")
$ws.Range("I3").Value = 46
$ws.Range("I6").Value = 64
$ws.Range("I7").Value = 159

$ws = $wb.Worksheets.Item("North Park")
$ws.Range("I3").Value = 14
$ws.Range("I7").Value = 48

$ws = $wb.Worksheets.Item("Morgan Park")
$ws.Range("I3").Value = 21
$ws.Range("I7").Value = 70

$ws = $wb.Worksheets.Item("Oakland")
$ws.Range("I2").Value = 14
$ws.Range("I7").Value = 45

$ws = $wb.Worksheets.Item("Riverdale")
$ws.Range("I3").Value = 29
$ws.Range("I7").Value = 83

$ws = $wb.Worksheets.Item("Auburn Gresham")
$ws.Range("I2").Value = 161
$ws.Range("I3").Value = 149
$ws.Range("I5").Value = 23
$ws.Range("I7").Value = 475
